$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("Goals of the project") - Content Placeholder 2, paragraph 1
#   "Ability to scan barcodes of new inventory from a cellphone camera"
#   -> "Ability to scan barcodes of inventory from a cellphone camera"
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$tr2.Paragraphs(1).Runs(1).Text = "Ability to scan barcodes of inventory from a cellphone camera"

# ---------------------------------------------------------------------------
# Slide 3 ("Interface") - Content Placeholder 5, paragraphs 1-3
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$shape3 = $slide3.Shapes.Item(5)
$tr3 = $shape3.TextFrame.TextRange
$tr3.Paragraphs(1).Runs(1).Text = "We were working on a test on mobile"
$tr3.Paragraphs(2).Runs(1).Text = "All features are easily accessible from the home page "
$tr3.Paragraphs(3).Runs(1).Text = "Interface complies with most basic und style guidelines "

# ---------------------------------------------------------------------------
# Slide 8 ("Scanner") - Content Placeholder 5, paragraph 1
#   Split the single run into three runs so the middle one ("Zxing") carries
#   its own run properties (with the spell-check err flag).
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$shape8 = $slide8.Shapes.Item(5)
$tr8 = $shape8.TextFrame.TextRange
$para8 = $tr8.Paragraphs(1)
$run8a = $para8.Runs(1)
$run8a.Text = "On android scanning is handled by a call to the external app “"
$run8b = $run8a.InsertAfter("Zxing")
$run8c = $run8b.InsertAfter(" Barcode Scanner”")

# ---------------------------------------------------------------------------
# Slide 9 ("Possible future work") - Content Placeholder 2
#   Insert a new paragraph after "Further optimizing the interface for
#   mobile " and before "Converting the data from the old system to this one"
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$shape9 = $slide9.Shapes.Item(2)
$tr9 = $shape9.TextFrame.TextRange
$para9b = $tr9.Paragraphs(2)
$run9b = $para9b.Runs(1)
$newPara9 = $run9b.InsertAfter("`rAdjusting the interface to follow more UND standards")
